$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") date values from 45750 to 45752 for rows 2 through 43
$ws.Range("C2:C43").Value = 45752
